$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert two new rows (4 and 5) for the new samples "Holden" and
#    "Rizzie Spiral". This shifts the old rows 4..29 down to 6..31.
# ------------------------------------------------------------------
$ws.Rows("4:5").Insert()

# Re-apply the same cell style used by the rest of column A (and used
# generally across the sheet) to the newly inserted row's cells so we
# don't leave cells with a different/blank style than their neighbours.
$ws.Range("A3").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)
$ws.Range("A3").Copy()
$ws.Range("B4:T5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2) Populate row 4 ("Holden", sample index 2)
# ------------------------------------------------------------------
$ws.Range("A4").Value2 = 2
$ws.Range("B4").Value2 = "Holden"
$ws.Range("C4").Value2 = 1.045371384473902
$ws.Range("D4").Value2 = 1.183026994165351
$ws.Range("E4").Value2 = 0.8637814456087102
$ws.Range("F4").Value2 = 1.183026994165351
$ws.Range("G4").Value2 = 0.8637814456087102
$ws.Range("H4").Value2 = 1.152382054292172
$ws.Range("I4").Value2 = 1.008881446077629
$ws.Range("J4").Value2 = 0.9727652774681128
$ws.Range("K4").Value2 = 0.8637814456087102
$ws.Range("L4").Value2 = 1.045371384473902
$ws.Range("M4").Value2 = 1.114199189319627
$ws.Range("N4").Value2 = 1.114199189319627
$ws.Range("O4").Value2 = 1.079093274905627
$ws.Range("P4").Value2 = 1.030726608082654
$ws.Range("Q4").Value2 = 1.030726608082654
$ws.Range("R4").Value2 = 0.9889903174641684
$ws.Range("S4").Value2 = 0.9889903174641684
$ws.Range("T4").Value2 = 1.03770143368098

# ------------------------------------------------------------------
# 3) Populate row 5 ("Rizzie Spiral", sample index 3)
# ------------------------------------------------------------------
$ws.Range("A5").Value2 = 3
$ws.Range("B5").Value2 = "Rizzie Spiral"
$ws.Range("C5").Value2 = 1.069598930381316
$ws.Range("D5").Value2 = 0.2458146743108857
$ws.Range("E5").Value2 = 1.801357195120993
$ws.Range("F5").Value2 = 0.2458146743108857
$ws.Range("G5").Value2 = 1.801357195120993
$ws.Range("H5").Value2 = 1.126489434276462
$ws.Range("I5").Value2 = 0.478761888460251
$ws.Range("J5").Value2 = 1.207742982447128
$ws.Range("K5").Value2 = 1.801357195120993
$ws.Range("L5").Value2 = 1.069598930381316
$ws.Range("M5").Value2 = 0.657706802346101
$ws.Range("N5").Value2 = 0.657706802346101
$ws.Range("O5").Value2 = 0.5980584977174844
$ws.Range("P5").Value2 = 1.038923599937732
$ws.Range("Q5").Value2 = 1.038923599937731
$ws.Range("R5").Value2 = 1.229531998733547
$ws.Range("S5").Value2 = 1.229531998733547
$ws.Range("T5").Value2 = 0.9882941841661724

# ------------------------------------------------------------------
# 4) Rename sample "Thomas Hex" -> "Matthies Hex". After the row
#    insert above, that sample (originally row 9) now lives on row 11.
# ------------------------------------------------------------------
$ws.Range("B11").Value2 = "Matthies Hex"
